# Apply changes described by the commit:
# "added assets and queues. Transactions are now made by retrieving queue items
#  from queue on Orchestrator."
#
# Concretely, in the workbook this means:
#  - Orders sheet: sort the product/quantity table (A2:B8) alphabetically by product
#  - Orders sheet: record queue transaction results (Status/Notes) for each processed item
#  - Address sheet: phone number value is corrected/trimmed
#  - Orders tab becomes the active/selected sheet (instead of Address)

$wb = $excel.ActiveWorkbook

$orders = $wb.Worksheets.Item("Orders")
$address = $wb.Worksheets.Item("Address")

# Sort the order rows (Product, Quantity) alphabetically by Product name.
# Row 9 (Ipoh Coff) is outside of this sorted block and stays untouched.
$sortRange = $orders.Range("A2:B8")
$keyRange = $orders.Range("A2:A8")
$sortRange.Sort($keyRange, 1)

# Record that each queue item / transaction was processed successfully.
$orders.Range("C2:C8").Value = "#512568"
$orders.Range("D2:D8").Value = "Sucess"

# Fix up the phone number on the Address sheet.
$address.Range("F2").Value = "345-8769"

# Make the Orders sheet the active tab (it was Address before).
$orders.Activate()
